$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1 (Documentation/Examples.docx body): the "Collision:" bullet used to
# read "Demonstrates the Lightning collision engine." -- extend it to also
# mention the AABB class: "Demonstrates the Lightning collision engine and
# AABB class."
# ---------------------------------------------------------------------------
$collision = $d.Content
$collision.Find.Execute("Demonstrates the Lightning collision engine.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$collision.Collapse(0)
$collision.MoveStart(1, -1)
$collision.InsertBefore(" and AABB class")

# ---------------------------------------------------------------------------
# Edit 2 (header): drop the red "(Release Candidate 0 - Pre-release)"
# qualifier that used to follow the version number, so the header now simply
# reads "Version 1.1.0 " before the line break.
# ---------------------------------------------------------------------------
$hdrRange = $d.Sections(1).Headers(1).Range
$hdrRange.Find.Execute("(Release Candidate 0 - Pre-release)", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------------
# Edit 3 (header): the release date used to say just "October, 2022"; add
# the day of month so it reads "October 29, 2022".
# ---------------------------------------------------------------------------
$dateRange = $d.Sections(1).Headers(1).Range
$dateRange.Find.Execute("October,", `
    $false, $false, $false, $false, $false, $true, 1, $false, "October 29,", 2)
